# Add a "Pending" column (M) to the Fragments sheet, and a new
# "pending from db" mock-fragment row (row 16), matching the new
# "f16" test fragment that exercises dependency-pending lookups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the Pending column.
$ws.Range("M1").Value = "Pending"

# New data row describing fragment f16 / "pending from db".
$ws.Range("A16").Value = "f16"
$ws.Range("D16").Value = "pending from db"
$ws.Range("M16").Value = "y"

# Move the active selection the way the authored workbook shows it
# (bottom pane selection moved down to the new empty row below the data).
[void]$ws.Range("A17").Select()
